{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 1: insert \"What did we find from Classifications?\" (plus a\n// following blank paragraph) right before the \"3)  Classifications\"\n// Heading 2 paragraph.\n// ---------------------------------------------------------------------\nlet classificationsHeadingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Classifications\") !== -1 && /^\\s*3\\)\\s*Classifications\\s*$/.test(t.replace(/\\s+/g, \" \"))) {\n    classificationsHeadingIndex = i;\n    break;\n  }\n}\nif (classificationsHeadingIndex === -1) {\n  throw new Error(\"Could not find the '3)  Classifications' heading paragraph.\");\n}\n\n// The paragraph immediately before the heading is an existing blank\n// paragraph; inserting *after* it (rather than before the heading)\n// keeps the new paragraph's style as the default/Normal style instead\n// of inheriting the Heading 2 style from the following paragraph.\nconst blankBeforeHeading = paragraphs.items[classificationsHeadingIndex - 1];\nconst newTextPara = blankBeforeHeading.insertParagraph(\n  \"What did we find from Classifications?\",\n  Word.InsertLocation.after\n);\nnewTextPara.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 2: replace the placeholder \"? \u2013 sent email to Henry asking\n// what this means.\" paragraph with the Perth/Brisbane IT job market\n// write-up, and remove the blank paragraph that used to follow it.\n// ---------------------------------------------------------------------\nconst body2 = context.document.body;\nconst paragraphs2 = body2.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet henryIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"sent email to Henry\") !== -1) {\n    henryIndex = i;\n    break;\n  }\n}\nif (henryIndex === -1) {\n  throw new Error(\"Could not find the 'sent email to Henry' paragraph.\");\n}\n\nconst henryParagraph = paragraphs2.items[henryIndex];\nconst blankAfterHenry = paragraphs2.items[henryIndex + 1];\n\n// Clear the placeholder sentence and type the new paragraph content.\nhenryParagraph.clear();\nawait context.sync();\n\n// These fragments mirror the individual runs introduced upstream; the\n// Word engine coalesces adjacent same-formatted runs back into a single\n// run when the package is saved, so they are joined into one string\n// before being typed in (this also avoids spurious xml:space=\"preserve\"\n// markers that inserting them one-by-one would otherwise leave behind).\nconst sentences = [\n  \"It is believed that the IT job market in Perth should be expanded.\",\n  \" Perth\\u2019s population is approx. 2 million whereas \",\n  \"Brisbane has a population of only approx. 2.3 million \\u2013 this is a \",\n  \"difference of \",\n  \"15% in size. Nevertheless\",\n  \", it \",\n  \"was found that\",\n  \" Perth \",\n  \"had only\",\n  \" 1092 IT job listings whereas Brisbane has 2,627 listings\",\n  \" \\u2013\",\n  \" this is \",\n  \"58.5% more listings in Brisbane. \",\n  \"Thus, it \",\n  \"appears there is significant room for the IT industry to grow in Perth.\"\n];\nhenryParagraph.insertText(sentences.join(\"\"), Word.InsertLocation.end);\nawait context.sync();\n\n// The blank paragraph that used to sit right after the Henry paragraph\n// is folded away in the final document.\nblankAfterHenry.delete();\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: insert \"What did we find from Classifications?\" (plus a\n# following blank paragraph) right before the \"3)  Classifications\"\n# Heading 2 paragraph.\n# ---------------------------------------------------------------------\n$headingIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    $norm = ($p.Range.Text -replace '\\s+', ' ').Trim()\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $norm -eq \"3) Classifications\") {\n        $headingIdx = $i\n        break\n    }\n}\nif ($headingIdx -eq -1) {\n    throw \"Could not find the '3)  Classifications' heading paragraph.\"\n}\n\n# The paragraph right before the heading is an existing blank paragraph.\n# Inserting two new paragraphs after it (rather than before the heading)\n# keeps their style as the default/Normal style instead of inheriting\n# the Heading 2 style from the following paragraph.\n$blankBeforeHeading = $d.Paragraphs($headingIdx - 1)\n$blankBeforeHeading.Range.InsertParagraphAfter()\n$blankBeforeHeading.Range.InsertParagraphAfter()\n$newTextPara = $d.Paragraphs($headingIdx - 1 + 1)\n$newTextPara.Range.Text = \"What did we find from Classifications?\"\n\n# ---------------------------------------------------------------------\n# Change 2: replace the placeholder \"? - sent email to Henry asking\n# what this means.\" paragraph with the Perth/Brisbane IT job market\n# write-up, and remove the blank paragraph that used to follow it.\n# ---------------------------------------------------------------------\n$henryIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -like \"*sent email to Henry*\") {\n        $henryIdx = $i\n        break\n    }\n}\nif ($henryIdx -eq -1) {\n    throw \"Could not find the 'sent email to Henry' paragraph.\"\n}\n\n$henryPara = $d.Paragraphs($henryIdx)\n$rightSingleQuote = [char]0x2019\n$enDash = [char]0x2013\n$newSentence = \"It is believed that the IT job market in Perth should be expanded.\" `\n    + \" Perth\" + $rightSingleQuote + \"s population is approx. 2 million whereas \" `\n    + \"Brisbane has a population of only approx. 2.3 million \" + $enDash + \" this is a \" `\n    + \"difference of \" `\n    + \"15% in size. Nevertheless\" `\n    + \", it \" `\n    + \"was found that\" `\n    + \" Perth \" `\n    + \"had only\" `\n    + \" 1092 IT job listings whereas Brisbane has 2,627 listings\" `\n    + \" \" + $enDash `\n    + \" this is \" `\n    + \"58.5% more listings in Brisbane. \" `\n    + \"Thus, it \" `\n    + \"appears there is significant room for the IT industry to grow in Perth.\"\n$henryPara.Range.Text = $newSentence\n\n# The blank paragraph that used to sit right after the Henry paragraph\n# is folded away in the final document.\n$blankAfterHenry = $d.Paragraphs($henryIdx + 1)\n$blankAfterHenry.Range.Delete()\n"}
